# Apply "update training config and result" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# New cells that did not exist before need their formatting copied from a
# same-column neighbour before the value is written (otherwise the engine
# falls back to the <col> default style and creates a brand-new cellXf).
# ---------------------------------------------------------------------------
function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# Row 8 -------------------------------------------------------------------
Copy-Format "B7" "C8"   # Epoch column style (plain number, style 1)
Copy-Format "D2" "D8"   # Learning Rate column style (scientific, style 4)
Copy-Format "J5" "J8"   # Comment column style (style 1)
Copy-Format "H1" "K8"   # new "extended comment" column - header style (style 2)

$ws.Range("A8").Value = 44841
$ws.Range("B8").Value = "N/A"
$ws.Range("C8").Value = 200
$ws.Range("D8").Value = 0.1
$ws.Range("E8").Value = "[0-14]"
$ws.Range("F8").Value = 512
$ws.Range("G8").Value = "N/A"
$ws.Range("H8").Value = "Adam"
$ws.Range("I8").Value = "6.05973(Invalid)"
$ws.Range("J8").Value = "Use filter data function" + [char]10 + "Normalization"
$ws.Range("K8").Value = "Did not normalize the testing data"

# Row 9 ---------------------------------------------------------------------
Copy-Format "D2" "D9"   # Learning Rate column style (scientific, style 4)
Copy-Format "J5" "J9"   # Comment column style (style 1)

$ws.Range("A9").Value = 44842
$ws.Range("B9").Value = "N/A"
$ws.Range("C9").Value = 200
$ws.Range("D9").Value = 0.1
$ws.Range("E9").Value = "[0-14]"
$ws.Range("F9").Value = 512
$ws.Range("G9").Value = "N/A"
$ws.Range("H9").Value = "Adam"
$ws.Range("I9").Value = 8.14767
$ws.Range("J9").Value = "Use filter data function" + [char]10 + "Normalization"

# Row 10 ----------------------------------------------------------------
Copy-Format "D2" "D10"  # Learning Rate column style (scientific, style 4)

$ws.Range("A10").Value = 44843
$ws.Range("B10").Value = "N/A"
$ws.Range("C10").Value = 100
$ws.Range("D10").Value = 0.1
$ws.Range("E10").Value = "[1-9, 13, 14]"
$ws.Range("F10").Value = 512
$ws.Range("G10").Value = "N/A"
$ws.Range("H10").Value = "Adam"
$ws.Range("I10").Value = 3.92847

# Row 11 ----------------------------------------------------------------
Copy-Format "D2" "D11"  # Learning Rate column style (scientific, style 4)

$ws.Range("A11").Value = 44843
$ws.Range("B11").Value = "N/A"
$ws.Range("C11").Value = 100
$ws.Range("D11").Value = 0.1
$ws.Range("E11").Value = "[1-4, 6-9, 13, 14]"
$ws.Range("F11").Value = 512
$ws.Range("G11").Value = "N/A"
$ws.Range("H11").Value = "Adam"
$ws.Range("I11").Value = 3.79697

# ---------------------------------------------------------------------------
# Row heights - the workbook-wide default-height refresh that accompanied
# this edit (every row ends up with an explicit height).
# ---------------------------------------------------------------------------
for ($r = 1; $r -le 4; $r++) {
    $ws.Rows.Item($r).RowHeight = 12.75
}
$ws.Rows.Item(8).RowHeight = 38.25
$ws.Rows.Item(9).RowHeight = 38.25
$ws.Rows.Item(10).RowHeight = 12.75
$ws.Rows.Item(11).RowHeight = 25.5
for ($r = 12; $r -le 41; $r++) {
    $ws.Rows.Item($r).RowHeight = 12.75
}

# ---------------------------------------------------------------------------
# Selection matches where the author's cursor ended up.
# ---------------------------------------------------------------------------
$ws.Range("C10").Select() | Out-Null
